$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'67.174.71"
$ws.Range("E2").Value = "  -3.58%  "

$ws.Range("D3").Value = "'3.754.56"
$ws.Range("E3").Value = "  +0.63%  "

$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  +0.02%  "

$ws.Range("D5").Value = "'591.05"
$ws.Range("E5").Value = "  -3.71%  "

$ws.Range("D6").Value = "'171.67"
$ws.Range("E6").Value = "  -4.05%  "

$ws.Range("D7").Value = "'3.750.36"
$ws.Range("E7").Value = "  +0.49%  "

$ws.Range("E8").Value = "  +0.10%  "

$ws.Range("D9").Value = "'0.520"
$ws.Range("E9").Value = "  -1.79%  "

$ws.Range("E10").Value = "  -4.49%  "

$ws.Range("D11").Value = "'6.28"
$ws.Range("E11").Value = "  -4.33%  "

$ws.Range("E12").Value = "  -4.31%  "

$ws.Range("E13").Value = "  -5.18%  "

$ws.Range("E14").Value = "  -4.20%  "

$ws.Range("D15").Value = "'4.367.33"
$ws.Range("E15").Value = "  +0.69%  "

$ws.Range("D16").Value = "'3.744.48"
$ws.Range("E16").Value = "  +0.73%  "

$ws.Range("D17").Value = "'67.285.43"
$ws.Range("E17").Value = "  -3.38%  "

$ws.Range("E18").Value = "  -4.76%  "

$ws.Range("E19").Value = "  -5.45%  "

$ws.Range("D20").Value = "'16.06"
$ws.Range("E20").Value = "  -1.58%  "

$ws.Range("D21").Value = "'485.30"
$ws.Range("E21").Value = "  -3.24%  "

$ws.Range("D22").Value = "'9.08"
$ws.Range("E22").Value = "  -0.66%  "

$ws.Range("E23").Value = "  +0.00%  "

$ws.Range("D24").Value = "'83.89"
$ws.Range("E24").Value = "  -2.64%  "

$ws.Range("E25").Value = "  -9.38%  "

$ws.Range("D26").Value = "'0.0000137"
$ws.Range("E26").Value = "  +0.71%  "

$ws.Range("D27").Value = "'12.19"
$ws.Range("E27").Value = "  -5.59%  "

$ws.Range("E28").Value = "  -10.27%  "

$ws.Range("E29").Value = "  +0.06%  "

$ws.Range("E30").Value = "  -0.21%  "

$ws.Range("E31").Value = "  -2.67%  "

$ws.Range("D32").Value = "'32.09"
$ws.Range("E32").Value = "  +5.79%  "

$ws.Range("D33").Value = "'7.69"
$ws.Range("E33").Value = "  -4.55%  "

$ws.Range("E34").Value = "  -5.01%  "

$ws.Range("D35").Value = "'1.00"
$ws.Range("E35").Value = "  +0.37%  "

$ws.Range("D36").Value = "'0.999"
$ws.Range("E36").Value = "  -4.78%  "

$ws.Range("E37").Value = "  -1.71%  "

$ws.Range("D38").Value = "'5.70"
$ws.Range("E38").Value = "  -6.36%  "

$ws.Range("B39").Value = "Bittensor"
$ws.Range("C39").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D39").Value = "'451.46"
$ws.Range("E39").Value = "  +3.37%  "

$ws.Range("B40").Value = "TheGraph"
$ws.Range("C40").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D40").Value = "'0.322"
$ws.Range("E40").Value = "  -7.64%  "

$ws.Range("D41").Value = "'48.68"
$ws.Range("E41").Value = "  -1.90%  "

$ws.Range("E42").Value = "  -4.23%  "

$ws.Range("D43").Value = "'2.84"
$ws.Range("E43").Value = "  -6.13%  "

$ws.Range("B44").Value = "Cosmos"
$ws.Range("C44").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D44").Value = "'8.21"
$ws.Range("E44").Value = "  -4.13%  "

$ws.Range("B45").Value = "Arweave"
$ws.Range("C45").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D45").Value = "'41.19"
$ws.Range("E45").Value = "  -10.25%  "

$ws.Range("D46").Value = "'2.799.36"
$ws.Range("E46").Value = "  -5.23%  "

$ws.Range("B47").Value = "USDe"
$ws.Range("C47").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D47").Value = "'1.00"
$ws.Range("E47").Value = "  +0.04%  "

$ws.Range("B48").Value = "Monero"
$ws.Range("C48").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D48").Value = "'139.42"
$ws.Range("E48").Value = "  +0.02%  "

$ws.Range("D49").Value = "'0.0346"
$ws.Range("E49").Value = "  -4.01%  "

$ws.Range("D50").Value = "'25.80"
$ws.Range("E50").Value = "  -4.74%  "

$ws.Range("D51").Value = "'23.02"
$ws.Range("E51").Value = "  +8.33%  "

